$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.842.52"
$ws.Range("E2").Value = "  +0.54%  "

$ws.Range("D3").Value = "2.660.49"
$ws.Range("E3").Value = "  +0.49%  "

$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").Value = "'604.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.37%  "

$ws.Range("D6").Value = "'148.32"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.94%  "

$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.09%  "

$ws.Range("E8").Value = "  +0.38%  "

$ws.Range("E9").Value = "  +2.10%  "

$ws.Range("D10").Value = "'5.61"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.09%  "

$ws.Range("D11").Value = "'0.371"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.16%  "

$ws.Range("E12").Value = "  +0.18%  "

$ws.Range("D13").Value = "'27.74"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.32%  "

$ws.Range("D14").Value = "3.139.76"
$ws.Range("E14").Value = "  +0.54%  "

$ws.Range("D15").Value = "63.714.77"

$ws.Range("E16").Value = "  +0.96%  "

$ws.Range("D17").Value = "2.650.74"
$ws.Range("E17").Value = "  +1.49%  "

$ws.Range("D18").Value = "'11.50"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.35%  "

$ws.Range("D19").Value = "'4.58"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.25%  "

$ws.Range("D20").Value = "'343.59"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.02%  "

$ws.Range("D21").Value = "'7.02"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.90%  "

$ws.Range("E22").Value = "  -0.14%  "

$ws.Range("D23").Value = "'5.58"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.03%  "

$ws.Range("D24").Value = "'67.05"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.71%  "

$ws.Range("D25").Value = "'1.70"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.26%  "

$ws.Range("E26").Value = "  +8.41%  "

$ws.Range("E27").Value = "  +0.14%  "

$ws.Range("D28").Value = "'559.42"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.62%  "

$ws.Range("E29").Value = "  -0.01%  "

$ws.Range("E30").Value = "  -0.15%  "

$ws.Range("D31").Value = "'7.90"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.09%  "

$ws.Range("E32").Value = "  +3.76%  "

$ws.Range("E33").Value = "  -2.54%  "

$ws.Range("D34").Value = "0.0₃0825"
$ws.Range("E34").Value = "  +2.49%  "

$ws.Range("D35").Value = "'5.22"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +8.17%  "

$ws.Range("D36").Value = "'167.09"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.88%  "

$ws.Range("E37").Value = "  +1.41%  "

$ws.Range("D38").Value = "'1.00"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.08%  "

$ws.Range("E39").Value = "  +8.97%  "

$ws.Range("D40").Value = "'19.21"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.80%  "

$ws.Range("E41").Value = "  -0.03%  "

$ws.Range("D42").Value = "'168.75"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.30%  "

$ws.Range("D43").Value = "'3.81"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.15%  "

$ws.Range("D44").Value = "'22.45"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.44%  "

$ws.Range("D45").Value = "'0.0577"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.85%  "

$ws.Range("D46").Value = "'0.633"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.96%  "

$ws.Range("D47").Value = "'0.0249"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.97%  "

$ws.Range("D48").Value = "'0.0966"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.86%  "

$ws.Range("B49").Value = "dogwifhat"
$ws.Range("C49").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D49").Value = "'1.89"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +11.46%  "

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'18.90"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.31%  "

$ws.Range("E51").Value = "  +2.94%  "
